# Scheduled-runner style refresh of cached market-board figures
# (currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfit columns)
# across the per-job Leve profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 9413.674999999999
$ws.Range("I40").Value = 8121.3335
$ws.Range("J40").Value = 9967.536
$ws.Range("K40").Value = 8121.3335
$ws.Range("L40").Value = 9967.536
$ws.Range("M40").Value = -7946.3335
$ws.Range("N40").Value = -10317.536

$ws.Range("H69").Value = 8278.875
$ws.Range("J69").Value = 8278.875
$ws.Range("L69").Value = 24836.625
$ws.Range("N69").Value = -26584.625

$ws.Range("H70").Value = 1923.8077
$ws.Range("I70").Value = 935.4
$ws.Range("J70").Value = 2159.1428
$ws.Range("K70").Value = 2806.2
$ws.Range("L70").Value = 6477.428400000001
$ws.Range("M70").Value = -2536.2
$ws.Range("N70").Value = -7017.428400000001

$ws.Range("H72").Value = 8278.875
$ws.Range("J72").Value = 8278.875
$ws.Range("L72").Value = 74509.875
$ws.Range("N72").Value = -83245.875

$ws.Range("H73").Value = 1923.8077
$ws.Range("I73").Value = 935.4
$ws.Range("J73").Value = 2159.1428
$ws.Range("K73").Value = 2806.2
$ws.Range("L73").Value = 6477.428400000001
$ws.Range("M73").Value = -1870.2
$ws.Range("N73").Value = -8349.428400000001

$ws.Range("H99").Value = 2305.25
$ws.Range("J99").Value = 2856.2
$ws.Range("L99").Value = 8568.599999999999
$ws.Range("N99").Value = -11564.6

$ws.Range("H107").Value = 1444.3103
$ws.Range("I107").Value = 803.1111
$ws.Range("J107").Value = 2493.5454
$ws.Range("K107").Value = 803.1111
$ws.Range("L107").Value = 2493.5454
$ws.Range("M107").Value = 1116.8889
$ws.Range("N107").Value = -6333.5454

$ws.Range("H132").Value = 2401.1724
$ws.Range("I132").Value = 1793.6154
$ws.Range("K132").Value = 5380.8462
$ws.Range("M132").Value = -2850.8462


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 42317.6
$ws.Range("I2").Value = 27897
$ws.Range("K2").Value = 27897
$ws.Range("M2").Value = -27784

$ws.Range("H32").Value = 1097.8462
$ws.Range("I32").Value = 1091.76
$ws.Range("K32").Value = 1091.76
$ws.Range("M32").Value = -804.76

$ws.Range("H61").Value = 6084.6816
$ws.Range("I61").Value = 4148.0557
$ws.Range("K61").Value = 4148.0557
$ws.Range("M61").Value = -3936.0557

$ws.Range("H116").Value = 42317.6
$ws.Range("I116").Value = 27897
$ws.Range("K116").Value = 27897
$ws.Range("M116").Value = -25603

$ws.Range("H122").Value = 3019.077
$ws.Range("I122").Value = 2702.923
$ws.Range("J122").Value = 3651.3845
$ws.Range("K122").Value = 8108.768999999999
$ws.Range("L122").Value = 10954.1535
$ws.Range("M122").Value = -5658.768999999999
$ws.Range("N122").Value = -15854.1535

$ws.Range("H132").Value = 8819.477000000001
$ws.Range("I132").Value = 6148.875
$ws.Range("K132").Value = 18446.625
$ws.Range("M132").Value = -15916.625

$ws.Range("H136").Value = 6084.6816
$ws.Range("I136").Value = 4148.0557
$ws.Range("K136").Value = 12444.1671
$ws.Range("M136").Value = -9894.167099999999


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 42317.6
$ws.Range("I3").Value = 27897
$ws.Range("K3").Value = 27897
$ws.Range("M3").Value = -27783

$ws.Range("H80").Value = 476.94116
$ws.Range("I80").Value = 396.7143
$ws.Range("K80").Value = 396.7143
$ws.Range("M80").Value = 601.2857

$ws.Range("H83").Value = 476.94116
$ws.Range("I83").Value = 396.7143
$ws.Range("K83").Value = 1983.5715
$ws.Range("M83").Value = 3008.4285

$ws.Range("H86").Value = 3871.1875
$ws.Range("I86").Value = 2119.4167
$ws.Range("K86").Value = 2119.4167
$ws.Range("M86").Value = -996.4167000000002

$ws.Range("H89").Value = 3871.1875
$ws.Range("I89").Value = 2119.4167
$ws.Range("K89").Value = 10597.0835
$ws.Range("M89").Value = -4981.083500000001

$ws.Range("H99").Value = 2915.5715
$ws.Range("I99").Value = 2747.6956
$ws.Range("J99").Value = 3687.8
$ws.Range("K99").Value = 2747.6956
$ws.Range("L99").Value = 3687.8
$ws.Range("M99").Value = -1249.6956
$ws.Range("N99").Value = -6683.8

$ws.Range("H134").Value = 2626.8286
$ws.Range("I134").Value = 1312.1538
$ws.Range("K134").Value = 3936.4614
$ws.Range("M134").Value = -1401.4614


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3057.3809
$ws.Range("I99").Value = 3019.4
$ws.Range("J99").Value = 3091.9092
$ws.Range("K99").Value = 3019.4
$ws.Range("L99").Value = 3091.9092
$ws.Range("M99").Value = -1521.4
$ws.Range("N99").Value = -6087.9092

$ws.Range("H126").Value = 3057.3809
$ws.Range("I126").Value = 3019.4
$ws.Range("J126").Value = 3091.9092
$ws.Range("K126").Value = 9058.200000000001
$ws.Range("L126").Value = 9275.7276
$ws.Range("M126").Value = -6588.200000000001
$ws.Range("N126").Value = -14215.7276

$ws.Range("H132").Value = 3941.6296
$ws.Range("I132").Value = 2729.6667
$ws.Range("K132").Value = 8189.000100000001
$ws.Range("M132").Value = -5659.000100000001

$ws.Range("H134").Value = 6216.32
$ws.Range("I134").Value = 6113.409
$ws.Range("K134").Value = 18340.227
$ws.Range("M134").Value = -15805.227


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1134.5
$ws.Range("J34").Value = 1852.75
$ws.Range("L34").Value = 5558.25
$ws.Range("N34").Value = -5726.25

$ws.Range("H39").Value = 5862.375
$ws.Range("I39").Value = 499.5
$ws.Range("J39").Value = 7650
$ws.Range("K39").Value = 1498.5
$ws.Range("L39").Value = 22950
$ws.Range("M39").Value = -1204.5
$ws.Range("N39").Value = -23538

$ws.Range("H55").Value = 2351
$ws.Range("J55").Value = 2688.75
$ws.Range("L55").Value = 8066.25
$ws.Range("N55").Value = -8420.25

$ws.Range("H86").Value = 826
$ws.Range("J86").Value = 1449.6
$ws.Range("L86").Value = 4348.799999999999
$ws.Range("N86").Value = -6720.799999999999

$ws.Range("H89").Value = 826
$ws.Range("J89").Value = 1449.6
$ws.Range("L89").Value = 13046.4
$ws.Range("N89").Value = -24902.4

$ws.Range("H132").Value = 5794.8184
$ws.Range("J132").Value = 6130.375
$ws.Range("L132").Value = 55173.375
$ws.Range("N132").Value = -60233.375

$ws.Range("H137").Value = 1002637.2
$ws.Range("I137").Value = 1429120
$ws.Range("J137").Value = 7510.6665
$ws.Range("K137").Value = 4287360
$ws.Range("L137").Value = 22531.9995
$ws.Range("M137").Value = -4282260
$ws.Range("N137").Value = -32731.9995


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 19999.5
$ws.Range("J19").Value = 19999.5
$ws.Range("L19").Value = 19999.5
$ws.Range("N19").Value = -20575.5

$ws.Range("H102").Value = 3516
$ws.Range("I102").Value = 2124.75
$ws.Range("K102").Value = 2124.75
$ws.Range("M102").Value = -502.75

$ws.Range("H122").Value = 7736.1934
$ws.Range("I122").Value = 6862.8887
$ws.Range("K122").Value = 20588.6661
$ws.Range("M122").Value = -18138.6661

$ws.Range("H126").Value = 7469.5
$ws.Range("I126").Value = 3388.6667
$ws.Range("K126").Value = 10166.0001
$ws.Range("M126").Value = -7696.000100000001

$ws.Range("H132").Value = 6273.5
$ws.Range("I132").Value = 4151.9
$ws.Range("K132").Value = 12455.7
$ws.Range("M132").Value = -9925.699999999999


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 34952.9
$ws.Range("I22").Value = 57271.57
$ws.Range("K22").Value = 57271.57
$ws.Range("M22").Value = -56976.57

$ws.Range("H27").Value = 34952.9
$ws.Range("I27").Value = 57271.57
$ws.Range("K27").Value = 57271.57
$ws.Range("M27").Value = -57164.57

$ws.Range("H46").Value = 4337.1
$ws.Range("I46").Value = 1474
$ws.Range("J46").Value = 7200.2
$ws.Range("K46").Value = 1474
$ws.Range("L46").Value = 7200.2
$ws.Range("M46").Value = -1286
$ws.Range("N46").Value = -7576.2

$ws.Range("H122").Value = 8889.9
$ws.Range("J122").Value = 11499.667
$ws.Range("L122").Value = 34499.001
$ws.Range("N122").Value = -39399.001

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0

$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 6433.021
$ws.Range("I13").Value = 5470.5884
$ws.Range("J13").Value = 8770.357
$ws.Range("K13").Value = 5470.5884
$ws.Range("L13").Value = 8770.357
$ws.Range("M13").Value = -5330.5884
$ws.Range("N13").Value = -9050.357

$ws.Range("H62").Value = 7874.875
$ws.Range("I62").Value = 7500
$ws.Range("J62").Value = 7999.8335
$ws.Range("K62").Value = 7500
$ws.Range("L62").Value = 7999.8335
$ws.Range("M62").Value = -6876
$ws.Range("N62").Value = -9247.833500000001

$ws.Range("H65").Value = 7874.875
$ws.Range("I65").Value = 7500
$ws.Range("J65").Value = 7999.8335
$ws.Range("K65").Value = 37500
$ws.Range("L65").Value = 39999.1675
$ws.Range("M65").Value = -34380
$ws.Range("N65").Value = -46239.1675

$ws.Range("H81").Value = 4328.909
$ws.Range("I81").Value = 2861.6
$ws.Range("K81").Value = 5723.2
$ws.Range("M81").Value = -4662.2

$ws.Range("H84").Value = 4328.909
$ws.Range("I84").Value = 2861.6
$ws.Range("K84").Value = 28616
$ws.Range("M84").Value = -23312

$ws.Range("H136").Value = 3117.5107
$ws.Range("I136").Value = 2779.325
$ws.Range("K136").Value = 8337.974999999999
$ws.Range("M136").Value = -5787.974999999999

